{"js": "// Review fix (\"Resultat af g\u00e5rdagens reviews\"): the cross-reference line read\n// \"L \u00e6s egen historik\" because a stray space had split the word \"L\u00e6s\" in two.\n// Remove that space so it reads \"L\u00e6s egen historik\" again.\nconst body = context.document.body;\n\nconst hits = body.search(\"L \u00e6s\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[0].insertText(\"L\u00e6s\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Review fix (\"Resultat af g\u00e5rdagens reviews\"): the cross-reference line read\n# \"L \u00e6s egen historik\" because a stray space had split the word \"L\u00e6s\" in two.\n# Remove that space so it reads \"L\u00e6s egen historik\" again.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"L \u00e6s\"\n$find.Replacement.Text = \"L\u00e6s\"\n$find.Execute([ref]$find.Text, $false, $true, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n"}
